# toArray added to excelWorkbook
# This script adds a handful of new "tag" cells scattered through Sheet1
# using the existing shared strings (fdsa/fds) plus five brand-new ones
# (dsa/sa/ds/saf/af) so the last-used row of the sheet extends further
# down (through row 88), letting the last-row detection be exercised.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Prime the shared-string table with the five new values in the order
# they were first typed during the original editing session.
$ws.Range("K6").Value = "dsa"
$ws.Range("L28").Value = "sa"
$ws.Range("J28").Value = "ds"
$ws.Range("F60").Value = "saf"
$ws.Range("E55").Value = "af"

# Row 2
$ws.Range("B2").Value = "fdsa"
$ws.Range("E2").Value = "fdsa"
$ws.Range("G2").Value = "fds"

# Row 3
$ws.Range("F3").Value = "fds"

# Row 4
$ws.Range("I4").Value = "fds"

# Row 5
$ws.Range("B5").Value = "fdsa"
$ws.Range("F5").Value = "fds"
$ws.Range("M5").Value = "fdsa"

# Row 6 (new row)
$ws.Range("G6").Value = "fds"
$ws.Range("L6").Value = "fds"

# Row 11
$ws.Range("L11").Value = "fdsa"

# Row 13
$ws.Range("L13").Value = "fdsa"

# Row 20
$ws.Range("K20").Value = "fds"
$ws.Range("L20").Value = "fds"

# Row 22
$ws.Range("J22").Value = "fdsa"

# Row 23
$ws.Range("P23").Value = "dsa"

# Row 28 (values already set above for shared-string priming)

# Row 31
$ws.Range("J31").Value = "dsa"

# Row 33
$ws.Range("O33").Value = "dsa"

# Row 34
$ws.Range("I34").Value = "dsa"

# Row 35
$ws.Range("L35").Value = "fds"

# Row 38
$ws.Range("M38").Value = "dsa"

# Row 41 (new row)
$ws.Range("J41").Value = "dsa"

# Row 46
$ws.Range("B46").Value = "fdsa"

# Row 48 (new row)
$ws.Range("D48").Value = "fdsa"
$ws.Range("J48").Value = "fdsa"

# Row 49 (new row)
$ws.Range("G49").Value = "dsa"

# Row 54 (new row)
$ws.Range("B54").Value = "sa"

# Row 55 (new row, E55 already set above)
$ws.Range("J55").Value = "dsa"

# Row 59 (new row)
$ws.Range("N59").Value = "fdsa"

# Row 60 (new row, F60 already set above)
$ws.Range("C60").Value = "ds"

# Row 61 (new row)
$ws.Range("B61").Value = "dsa"

# Row 66 (new row)
$ws.Range("G66").Value = "dsa"

# Row 72 (new row)
$ws.Range("P72").Value = "fdsa"

# Row 74 (new row)
$ws.Range("G74").Value = "fds"

# Row 83 (new row)
$ws.Range("Q83").Value = "fds"

# Row 84 (new row)
$ws.Range("F84").Value = "fdsa"

# Row 88 (new row)
$ws.Range("L88").Value = "dsa"

# Match the selection left behind by the original edit session.
[void]$ws.Range("N17").Select()

Write-Host "applied edits"
